$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '66.562.78'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '3.502.59'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '584.74'
$ws.Range('E5').Value = '  -2.22%  '
$ws.Range('D6').Value = '175.08'
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.496.42'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.596'
$ws.Range('E9').Value = '  -3.21%  '
$ws.Range('E10').Value = '  -4.70%  '
$ws.Range('D11').Value = '6.90'
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('E12').Value = '  -3.46%  '
$ws.Range('D13').Value = '4.107.27'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '30.29'
$ws.Range('E14').Value = '  -5.61%  '
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').Value = '66.552.73'
$ws.Range('E16').Value = '  -1.60%  '
$ws.Range('E17').Value = '  -3.47%  '
$ws.Range('D18').Value = '3.498.97'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').Value = '6.00'
$ws.Range('E19').Value = '  -5.12%  '
$ws.Range('E20').Value = '  -3.34%  '
$ws.Range('D21').Value = '380.18'
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('D22').Value = '7.86'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '72.12'
$ws.Range('E26').Value = '  -1.41%  '
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').Value = '  -4.98%  '
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = '24.40'
$ws.Range('E31').Value = '  +3.42%  '
$ws.Range('E32').Value = '  -4.43%  '
$ws.Range('D33').Value = '2.00'
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('E34').Value = '  -7.31%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').Value = '7.23'
$ws.Range('E36').Value = '  -2.61%  '
$ws.Range('E37').Value = '  -2.79%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '29.57'
$ws.Range('E38').Value = '  +12.51%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '159.95'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('D40').Value = '0.892'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('E41').Value = '  -5.69%  '
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('D43').Value = '6.44'
$ws.Range('E43').Value = '  -5.91%  '
$ws.Range('D44').Value = '2.54'
$ws.Range('E44').Value = '  -10.60%  '
$ws.Range('E45').Value = '  -4.42%  '
$ws.Range('D46').Value = '2.678.35'
$ws.Range('E46').Value = '  -5.41%  '
$ws.Range('D47').Value = '40.75'
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('D48').Value = '24.40'
$ws.Range('E48').Value = '  -8.99%  '
$ws.Range('D49').Value = '0.0291'
$ws.Range('E49').Value = '  -3.28%  '
$ws.Range('D50').Value = '315.27'
$ws.Range('E50').Value = '  -6.22%  '
$ws.Range('E51').Value = '  -4.81%  '
